$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 10.70913233333333
$ws.Range("H2").Value = 32.127397
$ws.Range("I2").Value = 0.007451729107954897
$ws.Range("J2").Value = 0.007451729107954897
$ws.Range("M2").Value = 19.21315233333334
$ws.Range("N2").Value = 57.63945700000001
$ws.Range("O2").Value = 0.04451179209991234
$ws.Range("P2").Value = 0.04451179209991233
$ws.Range("Q2").Value = 205.7561908781588
$ws.Range("R2").Value = 1851.805717903429
$ws.Range("S2").Value = 0.0003316898168381536
$ws.Range("T2").Value = 0.0003316898168381535
$ws.Range("G3").Value = 10.70913233333333
$ws.Range("H3").Value = 32.127397
$ws.Range("I3").Value = 0.007451729107954897
$ws.Range("J3").Value = 0.007451729107954897
$ws.Range("O3").Value = 0.2141755495962477
$ws.Range("P3").Value = 0.2141755495962477
$ws.Range("Q3").Value = 990.0285561462916
$ws.Range("R3").Value = 8910.257005316624
$ws.Range("S3").Value = 0.001595978177138597
$ws.Range("T3").Value = 0.001595978177138597
$ws.Range("G4").Value = 10.70913233333333
$ws.Range("H4").Value = 32.127397
$ws.Range("I4").Value = 0.007451729107954897
$ws.Range("J4").Value = 0.007451729107954897
$ws.Range("M4").Value = 166.8580016666666
$ws.Range("N4").Value = 500.5740049999999
$ws.Range("O4").Value = 0.3865658561145097
$ws.Range("P4").Value = 0.3865658561145097
$ws.Range("Q4").Value = 1786.904420723887
$ws.Range("R4").Value = 16082.13978651498
$ws.Range("S4").Value = 0.002880584042149996
$ws.Range("T4").Value = 0.002880584042149996
$ws.Range("G5").Value = 10.70913233333333
$ws.Range("H5").Value = 32.127397
$ws.Range("I5").Value = 0.007451729107954897
$ws.Range("J5").Value = 0.007451729107954897
$ws.Range("M5").Value = 41.09915599999999
$ws.Range("N5").Value = 123.297468
$ws.Range("O5").Value = 0.09521587377309249
$ws.Range("P5").Value = 0.09521587377309249
$ws.Range("Q5").Value = 440.1363003923107
$ws.Range("R5").Value = 3961.226703530795
$ws.Range("S5").Value = 0.0007095228981343124
$ws.Range("T5").Value = 0.0007095228981343124
$ws.Range("G6").Value = 10.70913233333333
$ws.Range("H6").Value = 32.127397
$ws.Range("I6").Value = 0.007451729107954897
$ws.Range("J6").Value = 0.007451729107954897
$ws.Range("M6").Value = 112.0244103333333
$ws.Range("N6").Value = 336.073231
$ws.Range("O6").Value = 0.2595309284162377
$ws.Range("P6").Value = 0.2595309284162377
$ws.Range("Q6").Value = 1199.684234823301
$ws.Range("R6").Value = 10797.15811340971
$ws.Range("S6").Value = 0.001933954173693837
$ws.Range("T6").Value = 0.001933954173693837
$ws.Range("I7").Value = 0.03290895798513831
$ws.Range("J7").Value = 0.03290895798513832
$ws.Range("M7").Value = 19.21315233333334
$ws.Range("N7").Value = 57.63945700000001
$ws.Range("O7").Value = 0.04451179209991234
$ws.Range("P7").Value = 0.04451179209991233
$ws.Range("Q7").Value = 908.6779380590996
$ws.Range("R7").Value = 8178.101442531896
$ws.Range("S7").Value = 0.001464836696059227
$ws.Range("T7").Value = 0.001464836696059227
$ws.Range("I8").Value = 0.03290895798513831
$ws.Range("J8").Value = 0.03290895798513832
$ws.Range("O8").Value = 0.2141755495962477
$ws.Range("P8").Value = 0.2141755495962477
$ws.Range("S8").Value = 0.007048294163106824
$ws.Range("T8").Value = 0.007048294163106824
$ws.Range("I9").Value = 0.03290895798513831
$ws.Range("J9").Value = 0.03290895798513832
$ws.Range("M9").Value = 166.8580016666666
$ws.Range("N9").Value = 500.5740049999999
$ws.Range("O9").Value = 0.3865658561145097
$ws.Range("P9").Value = 0.3865658561145097
$ws.Range("Q9").Value = 7891.47883036763
$ws.Range("R9").Value = 71023.30947330866
$ws.Range("S9").Value = 0.01272147951736142
$ws.Range("T9").Value = 0.01272147951736143
$ws.Range("I10").Value = 0.03290895798513831
$ws.Range("J10").Value = 0.03290895798513832
$ws.Range("M10").Value = 41.09915599999999
$ws.Range("N10").Value = 123.297468
$ws.Range("O10").Value = 0.09521587377309249
$ws.Range("P10").Value = 0.09521587377309249
$ws.Range("Q10").Value = 1943.767252875886
$ws.Range("R10").Value = 17493.90527588298
$ws.Range("S10").Value = 0.003133455189516934
$ws.Range("T10").Value = 0.003133455189516935
$ws.Range("I11").Value = 0.03290895798513831
$ws.Range("J11").Value = 0.03290895798513832
$ws.Range("M11").Value = 112.0244103333333
$ws.Range("N11").Value = 336.073231
$ws.Range("O11").Value = 0.2595309284162377
$ws.Range("P11").Value = 0.2595309284162377
$ws.Range("Q11").Value = 5298.14724975531
$ws.Range("R11").Value = 47683.32524779779
$ws.Range("S11").Value = 0.008540892419093907
$ws.Range("T11").Value = 0.008540892419093908
$ws.Range("G12").Value = 411.37678
$ws.Range("H12").Value = 1234.13034
$ws.Range("I12").Value = 0.2862480573072345
$ws.Range("J12").Value = 0.2862480573072345
$ws.Range("M12").Value = 19.21315233333334
$ws.Range("N12").Value = 57.63945700000001
$ws.Range("O12").Value = 0.04451179209991234
$ws.Range("P12").Value = 0.04451179209991233
$ws.Range("Q12").Value = 7903.844740536155
$ws.Range("R12").Value = 71134.60266482539
$ws.Range("S12").Value = 0.01274141401586341
$ws.Range("T12").Value = 0.01274141401586341
$ws.Range("G13").Value = 411.37678
$ws.Range("H13").Value = 1234.13034
$ws.Range("I13").Value = 0.2862480573072345
$ws.Range("J13").Value = 0.2862480573072345
$ws.Range("O13").Value = 0.2141755495962477
$ws.Range("P13").Value = 0.2141755495962477
$ws.Range("Q13").Value = 38030.60293389259
$ws.Range("R13").Value = 342275.4264050333
$ws.Range("S13").Value = 0.06130733499463516
$ws.Range("T13").Value = 0.06130733499463515
$ws.Range("G14").Value = 411.37678
$ws.Range("H14").Value = 1234.13034
$ws.Range("I14").Value = 0.2862480573072345
$ws.Range("J14").Value = 0.2862480573072345
$ws.Range("M14").Value = 166.8580016666666
$ws.Range("N14").Value = 500.5740049999999
$ws.Range("O14").Value = 0.3865658561145097
$ws.Range("P14").Value = 0.3865658561145097
$ws.Range("Q14").Value = 68641.50744286795
$ws.Range("R14").Value = 617773.5669858116
$ws.Range("S14").Value = 0.1106537253340863
$ws.Range("T14").Value = 0.1106537253340863
$ws.Range("G15").Value = 411.37678
$ws.Range("H15").Value = 1234.13034
$ws.Range("I15").Value = 0.2862480573072345
$ws.Range("J15").Value = 0.2862480573072345
$ws.Range("M15").Value = 41.09915599999999
$ws.Range("N15").Value = 123.297468
$ws.Range("O15").Value = 0.09521587377309249
$ws.Range("P15").Value = 0.09521587377309249
$ws.Range("Q15").Value = 16907.23845599768
$ws.Range("R15").Value = 152165.1461039791
$ws.Range("S15").Value = 0.02725535889235858
$ws.Range("T15").Value = 0.02725535889235858
$ws.Range("G16").Value = 411.37678
$ws.Range("H16").Value = 1234.13034
$ws.Range("I16").Value = 0.2862480573072345
$ws.Range("J16").Value = 0.2862480573072345
$ws.Range("M16").Value = 112.0244103333333
$ws.Range("N16").Value = 336.073231
$ws.Range("O16").Value = 0.2595309284162377
$ws.Range("P16").Value = 0.2595309284162377
$ws.Range("Q16").Value = 46084.24120432539
$ws.Range("R16").Value = 414758.1708389285
$ws.Range("S16").Value = 0.07429022407029097
$ws.Range("T16").Value = 0.07429022407029097
$ws.Range("G17").Value = 173.2560603333334
$ws.Range("H17").Value = 519.768181
$ws.Range("I17").Value = 0.12055666021578
$ws.Range("J17").Value = 0.12055666021578
$ws.Range("M17").Value = 19.21315233333334
$ws.Range("N17").Value = 57.63945700000001
$ws.Range("O17").Value = 0.04451179209991234
$ws.Range("P17").Value = 0.04451179209991233
$ws.Range("Q17").Value = 3328.795079857525
$ws.Range("R17").Value = 29959.15571871772
$ws.Range("S17").Value = 0.005366192995784571
$ws.Range("T17").Value = 0.005366192995784571
$ws.Range("G18").Value = 173.2560603333334
$ws.Range("H18").Value = 519.768181
$ws.Range("I18").Value = 0.12055666021578
$ws.Range("J18").Value = 0.12055666021578
$ws.Range("O18").Value = 0.2141755495962477
$ws.Range("P18").Value = 0.2141755495962477
$ws.Range("Q18").Value = 16017.02564842755
$ws.Range("R18").Value = 144153.230835848
$ws.Range("S18").Value = 0.02582028895920277
$ws.Range("T18").Value = 0.02582028895920277
$ws.Range("G19").Value = 173.2560603333334
$ws.Range("H19").Value = 519.768181
$ws.Range("I19").Value = 0.12055666021578
$ws.Range("J19").Value = 0.12055666021578
$ws.Range("M19").Value = 166.8580016666666
$ws.Range("N19").Value = 500.5740049999999
$ws.Range("O19").Value = 0.3865658561145097
$ws.Range("P19").Value = 0.3865658561145097
$ws.Range("Q19").Value = 28909.16000385943
$ws.Range("R19").Value = 260182.4400347349
$ws.Range("S19").Value = 0.04660308856661904
$ws.Range("T19").Value = 0.04660308856661904
$ws.Range("G20").Value = 173.2560603333334
$ws.Range("H20").Value = 519.768181
$ws.Range("I20").Value = 0.12055666021578
$ws.Range("J20").Value = 0.12055666021578
$ws.Range("M20").Value = 41.09915599999999
$ws.Range("N20").Value = 123.297468
$ws.Range("O20").Value = 0.09521587377309249
$ws.Range("P20").Value = 0.09521587377309249
$ws.Range("Q20").Value = 7120.677851585078
$ws.Range("R20").Value = 64086.1006642657
$ws.Range("S20").Value = 0.01147890774161131
$ws.Range("T20").Value = 0.01147890774161131
$ws.Range("G21").Value = 173.2560603333334
$ws.Range("H21").Value = 519.768181
$ws.Range("I21").Value = 0.12055666021578
$ws.Range("J21").Value = 0.12055666021578
$ws.Range("M21").Value = 112.0244103333333
$ws.Range("N21").Value = 336.073231
$ws.Range("O21").Value = 0.2595309284162377
$ws.Range("P21").Value = 0.2595309284162377
$ws.Range("Q21").Value = 19408.90799551809
$ws.Range("R21").Value = 174680.1719596628
$ws.Range("S21").Value = 0.03128818195256228
$ws.Range("T21").Value = 0.03128818195256229
$ws.Range("G22").Value = 794.4973246666667
$ws.Range("H22").Value = 2383.491974
$ws.Range("I22").Value = 0.5528345953838922
$ws.Range("J22").Value = 0.5528345953838923
$ws.Range("M22").Value = 19.21315233333334
$ws.Range("N22").Value = 57.63945700000001
$ws.Range("O22").Value = 0.04451179209991234
$ws.Range("P22").Value = 0.04451179209991233
$ws.Range("Q22").Value = 15264.79812724646
$ws.Range("R22").Value = 137383.1831452182
$ws.Range("S22").Value = 0.02460765857536697
$ws.Range("T22").Value = 0.02460765857536697
$ws.Range("G23").Value = 794.4973246666667
$ws.Range("H23").Value = 2383.491974
$ws.Range("I23").Value = 0.5528345953838922
$ws.Range("J23").Value = 0.5528345953838923
$ws.Range("O23").Value = 0.2141755495962477
$ws.Range("P23").Value = 0.2141755495962477
$ws.Range("Q23").Value = 73448.99798777643
$ws.Range("R23").Value = 661040.9818899878
$ws.Range("S23").Value = 0.1184036533021644
$ws.Range("T23").Value = 0.1184036533021644
$ws.Range("G24").Value = 794.4973246666667
$ws.Range("H24").Value = 2383.491974
$ws.Range("I24").Value = 0.5528345953838922
$ws.Range("J24").Value = 0.5528345953838923
$ws.Range("M24").Value = 166.8580016666666
$ws.Range("N24").Value = 500.5740049999999
$ws.Range("O24").Value = 0.3865658561145097
$ws.Range("P24").Value = 0.3865658561145097
$ws.Range("Q24").Value = 132568.2359233929
$ws.Range("R24").Value = 1193114.123310536
$ws.Range("S24").Value = 0.2137069786542929
$ws.Range("T24").Value = 0.2137069786542929
$ws.Range("G25").Value = 794.4973246666667
$ws.Range("H25").Value = 2383.491974
$ws.Range("I25").Value = 0.5528345953838922
$ws.Range("J25").Value = 0.5528345953838923
$ws.Range("M25").Value = 41.09915599999999
$ws.Range("N25").Value = 123.297468
$ws.Range("O25").Value = 0.09521587377309249
$ws.Range("P25").Value = 0.09521587377309249
$ws.Range("Q25").Value = 32653.16948805798
$ws.Range("R25").Value = 293878.5253925218
$ws.Range("S25").Value = 0.05263862905147134
$ws.Range("T25").Value = 0.05263862905147135
$ws.Range("G26").Value = 794.4973246666667
$ws.Range("H26").Value = 2383.491974
$ws.Range("I26").Value = 0.5528345953838922
$ws.Range("J26").Value = 0.5528345953838923
$ws.Range("M26").Value = 112.0244103333333
$ws.Range("N26").Value = 336.073231
$ws.Range("O26").Value = 0.2595309284162377
$ws.Range("P26").Value = 0.2595309284162377
$ws.Range("Q26").Value = 89003.09430719423
$ws.Range("R26").Value = 801027.8487647481
$ws.Range("S26").Value = 0.1434776758005967
$ws.Range("T26").Value = 0.1434776758005967
